$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row shared strings:
#    columns A..J (1..10) carry the "_old" suffix -> rename to "_FV2310"
#    columns L..U (12..21) carry the "_new" suffix -> rename to "_FV2404"
#    column K (11) is "diff" and stays unchanged.
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2310')
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2404')
}

# 2. Turn the used range A1:U93 into an Excel Table ("Table1") with an
#    autofilter, picking up the renamed headers for its column names.
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U93"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
